$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value2 = 163.5
$ws.Range("J9").Value2 = 200
$ws.Range("L9").Value2 = 200
$ws.Range("N9").Value2 = -538
$ws.Range("H17").Value2 = 1507.7693
$ws.Range("J17").Value2 = 1256.3158
$ws.Range("L17").Value2 = 3768.9474
$ws.Range("N17").Value2 = -4104.9474
$ws.Range("H33").Value2 = 142
$ws.Range("I33").Value2 = 202.5
$ws.Range("K33").Value2 = 202.5
$ws.Range("M33").Value2 = 26.5
$ws.Range("H38").Value2 = 471.44446
$ws.Range("I38").Value2 = 471.44446
$ws.Range("J38").Value2 = 0
$ws.Range("K38").Value2 = 1414.33338
$ws.Range("L38").Value2 = 0
$ws.Range("M38").Value2 = -1042.33338
$ws.Range("N38").Value2 = $null
$ws.Range("H43").Value2 = 1166.3334
$ws.Range("I43").Value2 = 0
$ws.Range("J43").Value2 = 1166.3334
$ws.Range("K43").Value2 = 0
$ws.Range("L43").Value2 = 1166.3334
$ws.Range("M43").Value2 = $null
$ws.Range("N43").Value2 = -1304.3334
$ws.Range("H61").Value2 = 0
$ws.Range("I61").Value2 = 0
$ws.Range("K61").Value2 = 0
$ws.Range("M61").Value2 = $null
$ws.Range("H62").Value2 = 1280
$ws.Range("I62").Value2 = 1280
$ws.Range("J62").Value2 = 0
$ws.Range("K62").Value2 = 1280
$ws.Range("L62").Value2 = 0
$ws.Range("M62").Value2 = -656
$ws.Range("N62").Value2 = $null
$ws.Range("H65").Value2 = 1280
$ws.Range("I65").Value2 = 1280
$ws.Range("J65").Value2 = 0
$ws.Range("K65").Value2 = 6400
$ws.Range("L65").Value2 = 0
$ws.Range("M65").Value2 = -3280
$ws.Range("N65").Value2 = -3280
$ws.Range("H74").Value2 = 3099.8
$ws.Range("I74").Value2 = 2749.5
$ws.Range("K74").Value2 = 2749.5
$ws.Range("M74").Value2 = -1813.5
$ws.Range("H77").Value2 = 3099.8
$ws.Range("I77").Value2 = 2749.5
$ws.Range("K77").Value2 = 13747.5
$ws.Range("M77").Value2 = -9067.5
$ws.Range("H100").Value2 = 1826.7273
$ws.Range("I100").Value2 = 1909.4
$ws.Range("K100").Value2 = 1909.4
$ws.Range("M100").Value2 = -1368.4
$ws.Range("H112").Value2 = 2409.55
$ws.Range("J112").Value2 = 2478.4736
$ws.Range("L112").Value2 = 7435.4208
$ws.Range("N112").Value2 = -9651.4208
$ws.Range("H116").Value2 = 13870.4
$ws.Range("I116").Value2 = 27625
$ws.Range("J116").Value2 = 4700.6665
$ws.Range("K116").Value2 = 27625
$ws.Range("L116").Value2 = 4700.6665
$ws.Range("M116").Value2 = -24183
$ws.Range("N116").Value2 = -11584.6665
$ws.Range("H127").Value2 = 2979.8
$ws.Range("I127").Value2 = 3570
$ws.Range("J127").Value2 = 2094.5
$ws.Range("K127").Value2 = 10710
$ws.Range("L127").Value2 = 6283.5
$ws.Range("M127").Value2 = -5750
$ws.Range("N127").Value2 = -16203.5
$ws.Range("H132").Value2 = 827.0612
$ws.Range("J132").Value2 = 1012.4286
$ws.Range("L132").Value2 = 3037.2858
$ws.Range("N132").Value2 = -8097.2858
$ws.Range("H135").Value2 = 349.08694
$ws.Range("I135").Value2 = 349.08694
$ws.Range("J135").Value2 = 0
$ws.Range("K135").Value2 = 3141.78246
$ws.Range("L135").Value2 = 0
$ws.Range("M135").Value2 = -606.7824600000004
$ws.Range("N135").Value2 = $null
$ws.Range("H137").Value2 = 1803.4814
$ws.Range("I137").Value2 = 1413.0834
$ws.Range("K137").Value2 = 4239.2502
$ws.Range("M137").Value2 = -1689.2502
$ws.Range("H138").Value2 = 2336.4238
$ws.Range("J138").Value2 = 2202.9092
$ws.Range("L138").Value2 = 6608.7276
$ws.Range("N138").Value2 = -16888.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value2 = 498.55554
$ws.Range("J5").Value2 = 899.75
$ws.Range("L5").Value2 = 899.75
$ws.Range("N5").Value2 = -1123.75
$ws.Range("H32").Value2 = 2420.747
$ws.Range("I32").Value2 = 1729.8334
$ws.Range("K32").Value2 = 1729.8334
$ws.Range("M32").Value2 = -1442.8334
$ws.Range("H61").Value2 = 14999
$ws.Range("I61").Value2 = 0
$ws.Range("J61").Value2 = 14999
$ws.Range("K61").Value2 = 0
$ws.Range("L61").Value2 = 14999
$ws.Range("M61").Value2 = $null
$ws.Range("N61").Value2 = -15423
$ws.Range("H74").Value2 = 1166
$ws.Range("I74").Value2 = 803.9474
$ws.Range("J74").Value2 = 1853.9
$ws.Range("K74").Value2 = 803.9474
$ws.Range("L74").Value2 = 1853.9
$ws.Range("M74").Value2 = 70.05259999999998
$ws.Range("N74").Value2 = -3601.9
$ws.Range("H77").Value2 = 1166
$ws.Range("I77").Value2 = 803.9474
$ws.Range("J77").Value2 = 1853.9
$ws.Range("K77").Value2 = 4019.737
$ws.Range("L77").Value2 = 9269.5
$ws.Range("M77").Value2 = 348.2629999999999
$ws.Range("N77").Value2 = -18005.5
$ws.Range("H102").Value2 = 1528.7693
$ws.Range("I102").Value2 = 1528.7693
$ws.Range("K102").Value2 = 1528.7693
$ws.Range("M102").Value2 = 93.23070000000007
$ws.Range("H122").Value2 = 1839.7778
$ws.Range("I122").Value2 = 1882.5
$ws.Range("J122").Value2 = 1498
$ws.Range("K122").Value2 = 5647.5
$ws.Range("L122").Value2 = 4494
$ws.Range("M122").Value2 = -3197.5
$ws.Range("N122").Value2 = -9394
$ws.Range("H132").Value2 = 1738.1951
$ws.Range("I132").Value2 = 1149.0667
$ws.Range("K132").Value2 = 3447.2001
$ws.Range("M132").Value2 = -917.2001
$ws.Range("H136").Value2 = 14999
$ws.Range("I136").Value2 = 0
$ws.Range("J136").Value2 = 14999
$ws.Range("K136").Value2 = 0
$ws.Range("L136").Value2 = 44997
$ws.Range("M136").Value2 = $null
$ws.Range("N136").Value2 = -50097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value2 = 498.55554
$ws.Range("J4").Value2 = 899.75
$ws.Range("L4").Value2 = 899.75
$ws.Range("N4").Value2 = -1129.75
$ws.Range("H20").Value2 = 1869.2307
$ws.Range("I20").Value2 = 1774.55
$ws.Range("K20").Value2 = 1774.55
$ws.Range("M20").Value2 = -1527.55
$ws.Range("H105").Value2 = 1870.4117
$ws.Range("I105").Value2 = 1870.4117
$ws.Range("K105").Value2 = 1870.4117
$ws.Range("M105").Value2 = -123.4117000000001
$ws.Range("H134").Value2 = 8316.652
$ws.Range("I134").Value2 = 11606.714
$ws.Range("J134").Value2 = 3198.7778
$ws.Range("K134").Value2 = 34820.142
$ws.Range("L134").Value2 = 9596.3334
$ws.Range("M134").Value2 = -32285.142
$ws.Range("N134").Value2 = -14666.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value2 = 0
$ws.Range("I6").Value2 = 0
$ws.Range("K6").Value2 = 0
$ws.Range("M6").Value2 = $null
$ws.Range("H7").Value2 = 102.9375
$ws.Range("I7").Value2 = 26.545454
$ws.Range("J7").Value2 = 271
$ws.Range("K7").Value2 = 26.545454
$ws.Range("L7").Value2 = 271
$ws.Range("M7").Value2 = 86.45454599999999
$ws.Range("N7").Value2 = -497
$ws.Range("H21").Value2 = 20000
$ws.Range("J21").Value2 = 20000
$ws.Range("L21").Value2 = 20000
$ws.Range("N21").Value2 = -20470
$ws.Range("H25").Value2 = 0
$ws.Range("I25").Value2 = 0
$ws.Range("K25").Value2 = 0
$ws.Range("M25").Value2 = $null
$ws.Range("H31").Value2 = 1983.75
$ws.Range("I31").Value2 = 1563.5264
$ws.Range("J31").Value2 = 2597.923
$ws.Range("K31").Value2 = 1563.5264
$ws.Range("L31").Value2 = 2597.923
$ws.Range("M31").Value2 = -1268.5264
$ws.Range("N31").Value2 = -3187.923
$ws.Range("H32").Value2 = 23000
$ws.Range("I32").Value2 = 0
$ws.Range("J32").Value2 = 23000
$ws.Range("K32").Value2 = 0
$ws.Range("L32").Value2 = 23000
$ws.Range("M32").Value2 = $null
$ws.Range("N32").Value2 = -23632
$ws.Range("H34").Value2 = 1983.75
$ws.Range("I34").Value2 = 1563.5264
$ws.Range("J34").Value2 = 2597.923
$ws.Range("K34").Value2 = 1563.5264
$ws.Range("L34").Value2 = 2597.923
$ws.Range("M34").Value2 = -1361.5264
$ws.Range("N34").Value2 = -3001.923
$ws.Range("H58").Value2 = 2175825
$ws.Range("I58").Value2 = 3106857
$ws.Range("K58").Value2 = 3106857
$ws.Range("M58").Value2 = -3106654
$ws.Range("H74").Value2 = 30999.6
$ws.Range("J74").Value2 = 30999.6
$ws.Range("L74").Value2 = 30999.6
$ws.Range("N74").Value2 = -32747.6
$ws.Range("H77").Value2 = 30999.6
$ws.Range("J77").Value2 = 30999.6
$ws.Range("L77").Value2 = 92998.79999999999
$ws.Range("N77").Value2 = -101734.8
$ws.Range("H105").Value2 = 3000
$ws.Range("I105").Value2 = 3000
$ws.Range("K105").Value2 = 3000
$ws.Range("M105").Value2 = -1253
$ws.Range("H122").Value2 = 2150.7273
$ws.Range("I122").Value2 = 920.2857
$ws.Range("K122").Value2 = 2760.8571
$ws.Range("M122").Value2 = -310.8571000000002
$ws.Range("H136").Value2 = 2175825
$ws.Range("I136").Value2 = 3106857
$ws.Range("K136").Value2 = 9320571
$ws.Range("M136").Value2 = -9318021

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value2 = 2500
$ws.Range("I81").Value2 = 2500
$ws.Range("K81").Value2 = 7500
$ws.Range("M81").Value2 = -6377
$ws.Range("H84").Value2 = 2500
$ws.Range("I84").Value2 = 2500
$ws.Range("K84").Value2 = 22500
$ws.Range("M84").Value2 = -16884
$ws.Range("H113").Value2 = 13388
$ws.Range("J113").Value2 = 1014.7143
$ws.Range("L113").Value2 = 3044.1429
$ws.Range("N113").Value2 = -7384.1429
$ws.Range("H131").Value2 = 10329.619
$ws.Range("J131").Value2 = 10938.329
$ws.Range("L131").Value2 = 32814.987
$ws.Range("N131").Value2 = -42894.987
$ws.Range("H141").Value2 = 5228.25
$ws.Range("I141").Value2 = 5228.25
$ws.Range("K141").Value2 = 15684.75
$ws.Range("M141").Value2 = -10504.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value2 = 2461363.2
$ws.Range("I126").Value2 = 6947182
$ws.Range("J126").Value2 = 68926.60000000001
$ws.Range("K126").Value2 = 20841546
$ws.Range("L126").Value2 = 206779.8
$ws.Range("M126").Value2 = -20839076
$ws.Range("N126").Value2 = -211719.8
$ws.Range("H132").Value2 = 1071018.1
$ws.Range("I132").Value2 = 1284425.2
$ws.Range("J132").Value2 = 3982.8333
$ws.Range("K132").Value2 = 3853275.6
$ws.Range("L132").Value2 = 11948.4999
$ws.Range("M132").Value2 = -3850745.6
$ws.Range("N132").Value2 = -17008.4999
$ws.Range("H139").Value2 = 59365.5
$ws.Range("J139").Value2 = 59365.5
$ws.Range("L139").Value2 = 59365.5
$ws.Range("N139").Value2 = -69645.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 2659.6843
$ws.Range("J7").Value2 = 9998.5
$ws.Range("L7").Value2 = 9998.5
$ws.Range("N7").Value2 = -10222.5
$ws.Range("H22").Value2 = 3527.2222
$ws.Range("I22").Value2 = 3949
$ws.Range("K22").Value2 = 3949
$ws.Range("M22").Value2 = -3654
$ws.Range("H27").Value2 = 3527.2222
$ws.Range("I27").Value2 = 3949
$ws.Range("K27").Value2 = 3949
$ws.Range("M27").Value2 = -3842
$ws.Range("H46").Value2 = 2171.0833
$ws.Range("J46").Value2 = 3120.6
$ws.Range("L46").Value2 = 3120.6
$ws.Range("N46").Value2 = -3496.6
$ws.Range("H68").Value2 = 3692.75
$ws.Range("I68").Value2 = 3257.1667
$ws.Range("K68").Value2 = 3257.1667
$ws.Range("M68").Value2 = -2508.1667
$ws.Range("H71").Value2 = 3692.75
$ws.Range("I71").Value2 = 3257.1667
$ws.Range("K71").Value2 = 16285.8335
$ws.Range("M71").Value2 = -12541.8335
$ws.Range("H126").Value2 = 2659.6843
$ws.Range("J126").Value2 = 9998.5
$ws.Range("L126").Value2 = 29995.5
$ws.Range("N126").Value2 = -34935.5
$ws.Range("H132").Value2 = 3078.5789
$ws.Range("J132").Value2 = 3516.5
$ws.Range("L132").Value2 = 10549.5
$ws.Range("N132").Value2 = -15609.5
$ws.Range("H136").Value2 = 4081
$ws.Range("I136").Value2 = 3040
$ws.Range("J136").Value2 = 5772.625
$ws.Range("K136").Value2 = 9120
$ws.Range("L136").Value2 = 17317.875
$ws.Range("M136").Value2 = -6570
$ws.Range("N136").Value2 = -22417.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value2 = 12000
$ws.Range("I32").Value2 = 12000
$ws.Range("K32").Value2 = 12000
$ws.Range("M32").Value2 = -11683
$ws.Range("H81").Value2 = 2712.75
$ws.Range("J81").Value2 = 2687.5
$ws.Range("L81").Value2 = 5375
$ws.Range("N81").Value2 = -7497
$ws.Range("H84").Value2 = 2712.75
$ws.Range("J84").Value2 = 2687.5
$ws.Range("L84").Value2 = 26875
$ws.Range("N84").Value2 = -37483
$ws.Range("H113").Value2 = 547.4286
$ws.Range("I113").Value2 = 322.46155
$ws.Range("J113").Value2 = 913
$ws.Range("K113").Value2 = 967.38465
$ws.Range("L113").Value2 = 2739
$ws.Range("M113").Value2 = 1202.61535
$ws.Range("N113").Value2 = -7079
$ws.Range("H119").Value2 = 24136.4
$ws.Range("J119").Value2 = 24136.4
$ws.Range("L119").Value2 = 24136.4
$ws.Range("N119").Value2 = -33812.4
$ws.Range("H122").Value2 = 61174.46
$ws.Range("I122").Value2 = 72033.91
$ws.Range("J122").Value2 = 1447.5
$ws.Range("K122").Value2 = 216101.73
$ws.Range("L122").Value2 = 4342.5
$ws.Range("M122").Value2 = -213651.73
$ws.Range("N122").Value2 = -9242.5
$ws.Range("H126").Value2 = 10795.23
$ws.Range("I126").Value2 = 12178.9
$ws.Range("K126").Value2 = 36536.7
$ws.Range("M126").Value2 = -34066.7
$ws.Range("H132").Value2 = 2574.889
$ws.Range("I132").Value2 = 1364.3334
$ws.Range("K132").Value2 = 4093.0002
$ws.Range("M132").Value2 = -1563.0002
$ws.Range("H136").Value2 = 42738584
$ws.Range("J136").Value2 = 2874.875
$ws.Range("L136").Value2 = 8624.625
$ws.Range("N136").Value2 = -13724.625
